$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: D1:P1 = سوال 4 .. سوال 16 ---
# Copy format from an existing header cell (A1) so the new header cells
# reuse the same bold/centered/bordered style used by A1:C1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1:P1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$headers = [ordered]@{
  "D1" = "سوال 4"
  "E1" = "سوال 5"
  "F1" = "سوال 6"
  "G1" = "سوال 7"
  "H1" = "سوال 8"
  "I1" = "سوال 9"
  "J1" = "سوال 10"
  "K1" = "سوال 11"
  "L1" = "سوال 12"
  "M1" = "سوال 13"
  "N1" = "سوال 14"
  "O1" = "سوال 15"
  "P1" = "سوال 16"
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# --- Row 11 ---
$row11 = [ordered]@{
  "B11" = "Option 2"
  "C11" = "Option 4"
  "D11" = "Option 3"
  "E11" = "Option 1"
  "F11" = "Option 3"
  "G11" = "Option 4"
  "H11" = "Option 3"
  "I11" = "Option 2"
  "J11" = "Option 1"
  "K11" = "Option 2"
  "L11" = "Option 4"
  "M11" = "Option 3"
  "N11" = "Option 2"
  "O11" = "Option 1"
  "P11" = "Option 4"
}
foreach ($addr in $row11.Keys) {
    $ws.Range($addr).Value = $row11[$addr]
}

# --- Row 12 ---
$row12 = [ordered]@{
  "B12" = "3-"
  "C12" = "4-"
  "D12" = "6-"
  "E12" = "2-"
  "F12" = "1-"
  "G12" = "4-"
  "H12" = "5-"
  "I12" = "6-"
  "J12" = "1-"
  "K12" = "2-"
  "L12" = "3-"
  "M12" = "6-"
  "N12" = "4-"
  "O12" = "5-"
  "P12" = "1-"
}
foreach ($addr in $row12.Keys) {
    $ws.Range($addr).Value = $row12[$addr]
}

# --- Row 13 ---
$row13 = [ordered]@{
  "B13" = "2-"
  "C13" = "6-"
  "D13" = "6-"
  "E13" = "5-"
  "F13" = "4-"
  "G13" = "5-"
  "H13" = "6-"
  "I13" = "5-"
  "J13" = "3-"
  "K13" = "6-"
  "L13" = "2-"
  "M13" = "1-"
  "N13" = "4-"
  "O13" = "5-"
  "P13" = "6-"
}
foreach ($addr in $row13.Keys) {
    $ws.Range($addr).Value = $row13[$addr]
}
